# PHOENIX-6078 UI : Added loop to insert the jurisdiction details in create employee feature
#
# The underlying code change added a loop that iterates the jurisdiction
# details when creating an employee. The accompanying functional-test data
# workbook (eisTestData.xlsx) was touched only cosmetically:
#   * the two header cells on the "jurisdictionList" sheet were
#     re-cased to camelCase ("JurisdictionType" -> "jurisdictionType",
#     "JurisdictionList" -> "jurisdictionList") so they line up with the
#     field names now produced by the loop, and
#   * the workbook was re-saved with the "jurisdictionList" sheet active
#     (instead of "assignmentDetails"), with a different cell selected on
#     each of those two sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Re-case the jurisdiction header cells on the jurisdictionList sheet
# ---------------------------------------------------------------------
$jurisdictionSheet = $wb.Worksheets.Item("jurisdictionList")
$jurisdictionSheet.Range("B1").Value = "jurisdictionType"
$jurisdictionSheet.Range("C1").Value = "jurisdictionList"

# ---------------------------------------------------------------------
# 2. Update the selected cell / active sheet bookkeeping to match the
#    state the workbook was saved in.
# ---------------------------------------------------------------------

# assignmentDetails: selection moved to H9, no longer the active tab
$assignmentSheet = $wb.Worksheets.Item("assignmentDetails")
$assignmentSheet.Range("H9").Select()

# jurisdictionList: selection moved to B18 and this sheet becomes active
$jurisdictionSheet.Range("B18").Select()
